$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Cells.Item(10, 9).Value = 'aa'
$ws.Cells.Item(10, 10).Value = 'Agree/Accept'
$ws.Cells.Item(11, 9).Value = 'aa'
$ws.Cells.Item(11, 10).Value = 'Agree/Accept'
$ws.Cells.Item(13, 9).Value = 'sv'
$ws.Cells.Item(13, 10).Value = 'Statement-opinion'
$ws.Cells.Item(15, 9).Value = 'aa'
$ws.Cells.Item(15, 10).Value = 'Agree/Accept'
$ws.Cells.Item(16, 9).Value = 'sd'
$ws.Cells.Item(16, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(24, 9).Value = 'sd'
$ws.Cells.Item(24, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(29, 9).Value = 'b'
$ws.Cells.Item(29, 10).Value = 'Acknowledge (Backchannel)'
$ws.Cells.Item(44, 9).Value = 'sd'
$ws.Cells.Item(44, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(46, 9).Value = 'sd'
$ws.Cells.Item(46, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(56, 9).Value = '%'
$ws.Cells.Item(56, 10).Value = 'Uninterpretable'
$ws.Cells.Item(59, 9).Value = 'sv'
$ws.Cells.Item(59, 10).Value = 'Statement-opinion'
$ws.Cells.Item(61, 9).Value = 'sd'
$ws.Cells.Item(61, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(75, 9).Value = 'b'
$ws.Cells.Item(75, 10).Value = 'Acknowledge (Backchannel)'
$ws.Cells.Item(80, 9).Value = 'sd'
$ws.Cells.Item(80, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(86, 9).Value = 'aa'
$ws.Cells.Item(86, 10).Value = 'Agree/Accept'
$ws.Cells.Item(95, 9).Value = 'aa'
$ws.Cells.Item(95, 10).Value = 'Agree/Accept'
$ws.Cells.Item(96, 9).Value = 'ba'
$ws.Cells.Item(96, 10).Value = 'Appreciation'
$ws.Cells.Item(101, 9).Value = 'sd'
$ws.Cells.Item(101, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(115, 9).Value = '%'
$ws.Cells.Item(115, 10).Value = 'Uninterpretable'
$ws.Cells.Item(116, 9).Value = 'aa'
$ws.Cells.Item(116, 10).Value = 'Agree/Accept'
$ws.Cells.Item(117, 9).Value = 'sv'
$ws.Cells.Item(117, 10).Value = 'Statement-opinion'
$ws.Cells.Item(121, 9).Value = 'aa'
$ws.Cells.Item(121, 10).Value = 'Agree/Accept'
$ws.Cells.Item(122, 9).Value = 'aa'
$ws.Cells.Item(122, 10).Value = 'Agree/Accept'
$ws.Cells.Item(123, 9).Value = 'aa'
$ws.Cells.Item(123, 10).Value = 'Agree/Accept'
$ws.Cells.Item(124, 9).Value = 'qy'
$ws.Cells.Item(124, 10).Value = 'Yes-No-Question'
$ws.Cells.Item(128, 9).Value = 'sd'
$ws.Cells.Item(128, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(142, 9).Value = 'sd'
$ws.Cells.Item(142, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(144, 9).Value = 'aa'
$ws.Cells.Item(144, 10).Value = 'Agree/Accept'
$ws.Cells.Item(148, 9).Value = 'aa'
$ws.Cells.Item(148, 10).Value = 'Agree/Accept'
$ws.Cells.Item(150, 9).Value = 'aa'
$ws.Cells.Item(150, 10).Value = 'Agree/Accept'
$ws.Cells.Item(152, 9).Value = 'aa'
$ws.Cells.Item(152, 10).Value = 'Agree/Accept'
$ws.Cells.Item(153, 9).Value = 'aa'
$ws.Cells.Item(153, 10).Value = 'Agree/Accept'
$ws.Cells.Item(155, 9).Value = 'aa'
$ws.Cells.Item(155, 10).Value = 'Agree/Accept'
$ws.Cells.Item(158, 9).Value = 'aa'
$ws.Cells.Item(158, 10).Value = 'Agree/Accept'
$ws.Cells.Item(166, 9).Value = 'aa'
$ws.Cells.Item(166, 10).Value = 'Agree/Accept'
$ws.Cells.Item(169, 9).Value = 'sd'
$ws.Cells.Item(169, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(172, 9).Value = '%'
$ws.Cells.Item(172, 10).Value = 'Uninterpretable'
$ws.Cells.Item(177, 9).Value = '%'
$ws.Cells.Item(177, 10).Value = 'Uninterpretable'
$ws.Cells.Item(184, 9).Value = '%'
$ws.Cells.Item(184, 10).Value = 'Uninterpretable'
$ws.Cells.Item(188, 9).Value = '%'
$ws.Cells.Item(188, 10).Value = 'Uninterpretable'
$ws.Cells.Item(190, 9).Value = '%'
$ws.Cells.Item(190, 10).Value = 'Uninterpretable'
$ws.Cells.Item(191, 9).Value = '%'
$ws.Cells.Item(191, 10).Value = 'Uninterpretable'
$ws.Cells.Item(192, 9).Value = '%'
$ws.Cells.Item(192, 10).Value = 'Uninterpretable'
$ws.Cells.Item(203, 9).Value = '%'
$ws.Cells.Item(203, 10).Value = 'Uninterpretable'
$ws.Cells.Item(213, 9).Value = 'aa'
$ws.Cells.Item(213, 10).Value = 'Agree/Accept'
$ws.Cells.Item(218, 9).Value = 'sd'
$ws.Cells.Item(218, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(221, 9).Value = 'sv'
$ws.Cells.Item(221, 10).Value = 'Statement-opinion'
$ws.Cells.Item(228, 9).Value = 'ba'
$ws.Cells.Item(228, 10).Value = 'Appreciation'
$ws.Cells.Item(230, 9).Value = 'sv'
$ws.Cells.Item(230, 10).Value = 'Statement-opinion'
$ws.Cells.Item(233, 9).Value = 'b'
$ws.Cells.Item(233, 10).Value = 'Acknowledge (Backchannel)'
$ws.Cells.Item(235, 9).Value = '%'
$ws.Cells.Item(235, 10).Value = 'Uninterpretable'
$ws.Cells.Item(240, 9).Value = 'aa'
$ws.Cells.Item(240, 10).Value = 'Agree/Accept'
$ws.Cells.Item(246, 9).Value = 'sd'
$ws.Cells.Item(246, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(247, 9).Value = 'sd'
$ws.Cells.Item(247, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(264, 9).Value = 'b'
$ws.Cells.Item(264, 10).Value = 'Acknowledge (Backchannel)'
$ws.Cells.Item(279, 9).Value = 'aa'
$ws.Cells.Item(279, 10).Value = 'Agree/Accept'
$ws.Cells.Item(282, 9).Value = 'ba'
$ws.Cells.Item(282, 10).Value = 'Appreciation'
$ws.Cells.Item(290, 9).Value = 'sd'
$ws.Cells.Item(290, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(335, 9).Value = 'aa'
$ws.Cells.Item(335, 10).Value = 'Agree/Accept'
$ws.Cells.Item(337, 9).Value = '%'
$ws.Cells.Item(337, 10).Value = 'Uninterpretable'
$ws.Cells.Item(345, 9).Value = 'sv'
$ws.Cells.Item(345, 10).Value = 'Statement-opinion'
